# Added Page Object and WebDriverKeyowrds
#
# - validLoginTest sheet: selection becomes A1:B2 (no longer the tab-selected sheet)
# - new "addEmployee" sheet added at the end, becomes the active/selected tab,
#   with Username/Password/Firstname/Lastname headers and a sample data row

$wb = $excel.ActiveWorkbook

# Update selection on the existing "validLoginTest" sheet and make sure it is
# no longer the tab-selected sheet (the new sheet will take that role).
$ws2 = $wb.Worksheets.Item("validLoginTest")
[void]$ws2.Range("A1:B2").Select()

# Append a brand new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "addEmployee"

# Header row
$ws3.Range("A1").Value = "Username"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "Firstname"
$ws3.Range("D1").Value = "Lastname"

# Data row
$ws3.Range("A2").Value = "Admin"
$ws3.Range("B2").Value = "admin123"
$ws3.Range("C2").Value = "test fname"
$ws3.Range("D2").Value = "test lname"

# Selection/active cell on the new sheet
[void]$ws3.Range("F7").Select()
